$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("YDS")

$ws.Range("B2").Value = "NCT(3.412677956402386, 1.3039830253046445, -0.6149790605244961, 2.738723303561531)"
$ws.Range("C2").Value = "JSU(-1.1034251727004434, 1.3400288549540393, 2.384430476837232, 6.346604452336397)"
$ws.Range("D2").Value = "NIG(1.2428863660165714, 0.9246837817812611, 1.2759474653387115, 2.8946451793456554)"
$ws.Range("E2").Value = "EXN(2.594712854689148, 2.5271088605615626, 2.85701988886492)"
